$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix classifications of stores: rows 207-257, column B (Classification), change 2 -> 3
$ws.Range("B207:B257").Value = 3

# Update the view state to match the saved selection/scroll position
$ws.Application.ActiveWindow.ScrollRow = 39
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A16").Select()
